$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the cells we will touch so that numeric-looking
# strings (e.g. "1.013") are stored as text, matching the source inlineStr type,
# then restore the "Normal" style so no stray formatting is left behind.
$touchedD = $ws.Range("D2:D51")
$touchedE = $ws.Range("E2:E51")
$touchedD.NumberFormat = "@"
$touchedE.NumberFormat = "@"

$ws.Range("D2").Value = "27.658.19"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.861.23"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.86%  "
$ws.Range("D5").Value = "332.70"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "0.4647"
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("D8").Value = "0.3887"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "46.28"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").Value = "0.07959"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "0.9946"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("D12").Value = "21.48"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").Value = "1.866.94"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "5.971"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "7.156"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "1.014"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "87.95"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "0.06704"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "0.00001042"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "16.87"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").Value = "27.654.25"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "5.440"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "10.85"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "2.318"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "158.41"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "19.66"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("D28").Value = "2.112"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "5.354"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").Value = "121.15"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").Value = "0.9686"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").Value = "0.09424"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "5.282"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "1.328"
$ws.Range("E35").Value = "  -8.45%  "
$ws.Range("D36").Value = "0.06008"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "1.195"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").Value = "8.125"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "1.010"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "0.5889"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").Value = "0.1876"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "0.5599"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("D46").Value = "12.09"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").Value = "1.910"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "3.287"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("D49").Value = "0.06767"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").Value = "111.88"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").Value = "1.061"
$ws.Range("E51").Value = "  -1.23%  "

$touchedD.Style = "Normal"
$touchedE.Style = "Normal"
